$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.709.97'
$ws.Range('E2').Value = '  +2.13%  '

$ws.Range('D3').Value = '2.387.25'
$ws.Range('E3').Value = '  +2.24%  '

$ws.Range('E4').Value = '  +0.05%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '552.35'
$ws.Range('E5').Value = '  +2.12%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.37'
$ws.Range('E6').Value = '  +4.58%  '

$ws.Range('E7').Value = '  -0.03%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.522'
$ws.Range('E8').Value = '  +0.01%  '

$ws.Range('D9').Value = '2.385.84'
$ws.Range('E9').Value = '  +2.18%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.108'
$ws.Range('E10').Value = '  +3.98%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.158'
$ws.Range('E11').Value = '  +1.95%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.36'
$ws.Range('E12').Value = '  +1.73%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.353'
$ws.Range('E13').Value = '  +3.99%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.91'
$ws.Range('E14').Value = '  +6.09%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000174'
$ws.Range('E15').Value = '  +8.66%  '

$ws.Range('D16').Value = '2.817.22'
$ws.Range('E16').Value = '  +2.22%  '

$ws.Range('D17').Value = '61.514.31'
$ws.Range('E17').Value = '  +1.73%  '

$ws.Range('D18').Value = '2.383.41'
$ws.Range('E18').Value = '  +2.07%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.11'
$ws.Range('E19').Value = '  +5.89%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.17'
$ws.Range('E20').Value = '  +2.96%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '322.90'
$ws.Range('E21').Value = '  +2.05%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.68'
$ws.Range('E22').Value = '  +2.14%  '

$ws.Range('E23').Value = '  +0.12%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '64.63'
$ws.Range('E24').Value = '  +2.76%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.75'
$ws.Range('E25').Value = '  -4.77%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.15'
$ws.Range('E26').Value = '  +6.76%  '

$ws.Range('B27').Value = 'Bittensor'
$ws.Range('C27').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '551.45'
$ws.Range('E27').Value = '  +10.64%  '

$ws.Range('B28').Value = 'Binance-PegBSC-USD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.998'
$ws.Range('E28').Value = '  +0.02%  '

$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.30'
$ws.Range('E29').Value = '  +4.84%  '

$ws.Range('B30').Value = 'WrappedeETH'
$ws.Range('C30').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D30').Value = '2.473.74'
$ws.Range('E30').Value = '  +0.96%  '

$ws.Range('D31').Value = '0.0₃0919'
$ws.Range('E31').Value = '  +3.46%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.42'
$ws.Range('E32').Value = '  +2.38%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.149'
$ws.Range('E33').Value = '  +3.00%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.85'
$ws.Range('E34').Value = '  +3.93%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.53'
$ws.Range('E35').Value = '  +1.12%  '

$ws.Range('B36').Value = 'RenderToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.76'
$ws.Range('E36').Value = '  +10.65%  '

$ws.Range('B37').Value = 'FirstDigitalUSD'
$ws.Range('C37').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.998'
$ws.Range('E37').Value = '  -0.13%  '

$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.94'
$ws.Range('E38').Value = '  +7.93%  '

$ws.Range('B39').Value = 'NEARProtocol'
$ws.Range('C39').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.74'
$ws.Range('E39').Value = '  +3.33%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.382'
$ws.Range('E40').Value = '  +2.74%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '18.56'
$ws.Range('E41').Value = '  +1.82%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '145.87'
$ws.Range('E42').Value = '  +6.12%  '

$ws.Range('E43').Value = '  +0.13%  '

$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.28'
$ws.Range('E44').Value = '  +8.66%  '

$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '148.86'
$ws.Range('E45').Value = '  +5.67%  '

$ws.Range('B46').Value = 'Filecoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.62'
$ws.Range('E46').Value = '  +2.74%  '

$ws.Range('B47').Value = 'Hedera'
$ws.Range('C47').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0530'
$ws.Range('E47').Value = '  +4.14%  '

$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '20.30'
$ws.Range('E48').Value = '  +5.01%  '

$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.585'
$ws.Range('E49').Value = '  +3.41%  '

$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0907'
$ws.Range('E50').Value = '  +1.19%  '

$ws.Range('B51').Value = 'VeChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0224'
$ws.Range('E51').Value = '  +1.77%  '
